# Applies the price/volume refresh captured in the commit "Updated cryptos
# list ... with GitHub Actions": column D (Price) and column E (Volume 1h)
# text cells are refreshed with newly scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text price strings (e.g. "37.846.51" / "59.31").
# Excel auto-detects simple decimals (single ".") as numbers, which would
# flip the cell from a text cell to a numeric one and introduce binary
# float noise (e.g. 59.31 -> 59.310000000000002). Forcing text format,
# writing the value, then clearing the format back keeps the cell a plain
# inline/shared string (matching the source workbook) without leaving any
# numbering/style attribute behind on the cell itself.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "D2" "37.846.51"
$ws.Range("E2").Value = "  +0.14%  "
Set-TextValue "D3" "2.082.47"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue "D5" "233.13"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("E6").Value = "  +0.23%  "
Set-TextValue "D7" "59.31"
$ws.Range("E7").Value = "  +3.39%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +2.02%  "
Set-TextValue "D10" "0.0788"
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("E12").Value = "  +2.39%  "
Set-TextValue "D13" "21.18"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("E15").Value = "  +2.61%  "
Set-TextValue "D16" "2.036.31"
$ws.Range("E16").Value = "  -2.15%  "
Set-TextValue "D17" "37.750.81"
$ws.Range("E17").Value = "  +0.25%  "
Set-TextValue "D18" "6.16"
$ws.Range("E18").Value = "  +0.34%  "
Set-TextValue "D19" "71.53"
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("E20").Value = "  +3.60%  "
Set-TextValue "D21" "228.13"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("E24").Value = "  +1.09%  "
Set-TextValue "D25" "171.56"
$ws.Range("E25").Value = "  +2.05%  "
Set-TextValue "D26" "9.15"
$ws.Range("E26").Value = "  +2.62%  "
Set-TextValue "D27" "0.137"
$ws.Range("E27").Value = "  -2.81%  "
$ws.Range("E28").Value = "  -1.27%  "
Set-TextValue "D29" "19.48"
$ws.Range("E29").Value = "  +0.14%  "
Set-TextValue "D31" "4.72"
$ws.Range("E31").Value = "  +2.52%  "
$ws.Range("E32").Value = "  +3.83%  "
$ws.Range("E33").Value = "  +1.24%  "
$ws.Range("E34").Value = "  +1.30%  "
Set-TextValue "D35" "3.43"
$ws.Range("E35").Value = "  +1.50%  "
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("E39").Value = "  -0.98%  "
Set-TextValue "D40" "99.04"
$ws.Range("E40").Value = "  +2.17%  "
$ws.Range("E41").Value = "  +1.79%  "
$ws.Range("E42").Value = "  -1.51%  "
Set-TextValue "D43" "16.69"
$ws.Range("E43").Value = "  +6.95%  "
Set-TextValue "D44" "1.444.02"
$ws.Range("E44").Value = "  -0.58%  "
Set-TextValue "D45" "1.15"
$ws.Range("E45").Value = "  -0.55%  "
Set-TextValue "D46" "4.18"
$ws.Range("E46").Value = "  +2.94%  "
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("E49").Value = "  -0.46%  "
Set-TextValue "D50" "2.274.33"
$ws.Range("E50").Value = "  -0.22%  "
Set-TextValue "D51" "46.84"
$ws.Range("E51").Value = "  +0.69%  "
